$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.761.64"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "'1.633.71"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'215.10"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'19.70"
$ws.Range("E10").Value = "  -3.53%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'4.26"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'1.858.47"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "'1.631.69"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "'62.80"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "'25.777.60"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'4.46"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "'194.46"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'9.95"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'6.29"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'1.82"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").Value = "'142.79"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "'15.57"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "'0.0495"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'3.34"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'0.903"
$ws.Range("D37").Value = "'1.128.67"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "'100.20"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "'1.767.62"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").Value = "'0.0₆0109"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").Value = "'55.12"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "'0.416"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'7.55"
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("E51").Value = "  +2.02%  "
